$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.968.33"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.642.58"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "212.74"
$ws.Range("D6").Value = "0.525"
$ws.Range("D8").Value = "23.53"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").Value = "1.874.96"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.646.49"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "27.958.32"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "233.48"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").Value = "153.32"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("D33").Value = "3.10"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "1.411.68"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "0.929"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "67.16"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "1.86"
$ws.Range("E44").Value = "  +7.06%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "1.784.51"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "88.01"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  -0.49%  "
